$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.999999981249027336538404142630
$ws.Range("C2").Value = 0.000000003055467619578928819814
$ws.Range("D2").Value = 0.000000002787152827383197825412
$ws.Range("E2").Value = 0.000000012908415468802579402637
$ws.Range("F2").Value = 46045
